$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.121.66"
$ws.Range("E2").Value = "  -3.43%  "
$ws.Range("D3").Value = "1.971.65"
$ws.Range("E3").Value = "  -3.21%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'243.00"
$ws.Range("E5").Value = "  -3.85%  "
$ws.Range("D6").Value = "'0.624"
$ws.Range("E6").Value = "  -3.21%  "
$ws.Range("D7").Value = "'61.59"
$ws.Range("E7").Value = "  -2.25%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").Value = "'0.372"
$ws.Range("E9").Value = "  -0.51%  "
$ws.Range("D10").Value = "'56.37"
$ws.Range("E10").Value = "  -4.35%  "
$ws.Range("D11").Value = "'0.0797"
$ws.Range("E11").Value = "  +5.82%  "
$ws.Range("E12").Value = "  -1.33%  "
$ws.Range("D13").Value = "'0.861"
$ws.Range("E13").Value = "  -5.21%  "
$ws.Range("D14").Value = "'21.97"
$ws.Range("E14").Value = "  +7.49%  "
$ws.Range("D15").Value = "'13.93"
$ws.Range("E15").Value = "  -7.11%  "
$ws.Range("D16").Value = "2.266.70"
$ws.Range("E16").Value = "  -3.05%  "
$ws.Range("D17").Value = "'5.41"
$ws.Range("E17").Value = "  -2.62%  "
$ws.Range("D18").Value = "1.986.24"
$ws.Range("E18").Value = "  -2.64%  "
$ws.Range("D19").Value = "35.966.31"
$ws.Range("E19").Value = "  -3.64%  "
$ws.Range("D20").Value = "'70.85"
$ws.Range("E20").Value = "  -3.13%  "
$ws.Range("D21").Value = "0.0₃0852"
$ws.Range("E21").Value = "  -2.16%  "
$ws.Range("D22").Value = "'238.21"
$ws.Range("E22").Value = "  +0.97%  "
$ws.Range("D23").Value = "'5.22"
$ws.Range("E23").Value = "  -1.71%  "
$ws.Range("D24").Value = "'1.00"
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").Value = "'2.49"
$ws.Range("E25").Value = "  -9.78%  "
$ws.Range("D26").Value = "'2.27"
$ws.Range("E26").Value = "  -2.41%  "
$ws.Range("D27").Value = "'9.72"
$ws.Range("E27").Value = "  +2.41%  "
$ws.Range("D28").Value = "'158.85"
$ws.Range("E28").Value = "  -3.99%  "
$ws.Range("D29").Value = "'0.134"
$ws.Range("E29").Value = "  +18.85%  "
$ws.Range("D30").Value = "'19.69"
$ws.Range("E30").Value = "  -0.41%  "
$ws.Range("E31").Value = "  -1.46%  "
$ws.Range("D32").Value = "'4.87"
$ws.Range("E32").Value = "  -5.81%  "
$ws.Range("E33").Value = "  -6.16%  "
$ws.Range("D34").Value = "'0.0616"
$ws.Range("E34").Value = "  +0.67%  "
$ws.Range("D35").Value = "'4.36"
$ws.Range("E35").Value = "  -6.78%  "
$ws.Range("D36").Value = "'6.24"
$ws.Range("E36").Value = "  +4.45%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("E38").Value = "  -6.08%  "
$ws.Range("D39").Value = "'1.83"
$ws.Range("E39").Value = "  +1.37%  "
$ws.Range("D40").Value = "'3.10"
$ws.Range("E40").Value = "  +14.71%  "
$ws.Range("D41").Value = "'0.0988"
$ws.Range("E41").Value = "  -5.29%  "
$ws.Range("D42").Value = "'1.22"
$ws.Range("B43").Value = "HuobiToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D43").Value = "'2.83"
$ws.Range("E43").Value = "  -3.51%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0212"
$ws.Range("E44").Value = "  -3.01%  "
$ws.Range("E45").Value = "  -4.17%  "
$ws.Range("D46").Value = "'92.67"
$ws.Range("E46").Value = "  -2.32%  "
$ws.Range("D47").Value = "'16.11"
$ws.Range("E47").Value = "  -4.06%  "
$ws.Range("D48").Value = "'7.51"
$ws.Range("E48").Value = "  -6.82%  "
$ws.Range("D49").Value = "1.346.13"
$ws.Range("E49").Value = "  -5.20%  "
$ws.Range("D50").Value = "'2.84"
$ws.Range("E50").Value = "  -3.22%  "
$ws.Range("D51").Value = "2.159.22"
$ws.Range("E51").Value = "  -3.06%  "
